$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the stats for row 24 (2025-11) per the latest data refresh
$ws.Range("B24").Value = 6357
$ws.Range("C24").Value = 1001
$ws.Range("D24").Value = 5959222
$ws.Range("E24").Value = 937.4267736353626
$ws.Range("F24").Value = 8.370269348789638
$ws.Range("G24").Value = 3.730569948186524
$ws.Range("H24").Value = 26.23940488567911
